$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 4.838602000000001
$ws.Range("H2").Value = 14.515806
$ws.Range("I2").Value = 0.04159701547454875
$ws.Range("J2").Value = 0.04159701547454875
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 29.57110033333333
$ws.Range("N2").Value = 88.713301
$ws.Range("O2").Value = 0.5311523066901919
$ws.Range("P2").Value = 0.5311523066901919
$ws.Range("Q2").Value = 143.0827852150674
$ws.Range("R2").Value = 1287.745066935606
$ws.Range("S2").Value = 0.02209435072073417
$ws.Range("T2").Value = 0.02209435072073417

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 4.838602000000001
$ws.Range("H3").Value = 14.515806
$ws.Range("I3").Value = 0.04159701547454875
$ws.Range("J3").Value = 0.04159701547454875
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 12.180337
$ws.Range("N3").Value = 36.541011
$ws.Range("O3").Value = 0.2187816490048282
$ws.Range("P3").Value = 0.2187816490048282
$ws.Range("Q3").Value = 58.93580296887401
$ws.Range("R3").Value = 530.4222267198661
$ws.Range("S3").Value = 0.009100663639201131
$ws.Range("T3").Value = 0.009100663639201131

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 4.838602000000001
$ws.Range("H4").Value = 14.515806
$ws.Range("I4").Value = 0.04159701547454875
$ws.Range("J4").Value = 0.04159701547454875
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 13.92204833333333
$ws.Range("N4").Value = 41.76614499999999
$ws.Range("O4").Value = 0.2500660443049799
$ws.Range("P4").Value = 0.2500660443049799
$ws.Range("Q4").Value = 67.36325090976334
$ws.Range("R4").Value = 606.2692581878699
$ws.Range("S4").Value = 0.01040200111461344
$ws.Range("T4").Value = 0.01040200111461344

$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 106.240761
$ws.Range("H5").Value = 318.722283
$ws.Range("I5").Value = 0.9133420313026024
$ws.Range("J5").Value = 0.9133420313026024
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 29.57110033333333
$ws.Range("N5").Value = 88.713301
$ws.Range("O5").Value = 0.5311523066901919
$ws.Range("P5").Value = 0.5311523066901919
$ws.Range("Q5").Value = 3141.656203020687
$ws.Range("R5").Value = 28274.90582718618
$ws.Range("S5").Value = 0.4851237267234827
$ws.Range("T5").Value = 0.4851237267234827

$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 106.240761
$ws.Range("H6").Value = 318.722283
$ws.Range("I6").Value = 0.9133420313026024
$ws.Range("J6").Value = 0.9133420313026024
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 12.180337
$ws.Range("N6").Value = 36.541011
$ws.Range("O6").Value = 0.2187816490048282
$ws.Range("P6").Value = 0.2187816490048282
$ws.Range("Q6").Value = 1294.048272116457
$ws.Range("R6").Value = 11646.43444904811
$ws.Range("S6").Value = 0.1998224757138028
$ws.Range("T6").Value = 0.1998224757138028

$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 106.240761
$ws.Range("H7").Value = 318.722283
$ws.Range("I7").Value = 0.9133420313026024
$ws.Range("J7").Value = 0.9133420313026024
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 13.92204833333333
$ws.Range("N7").Value = 41.76614499999999
$ws.Range("O7").Value = 0.2500660443049799
$ws.Range("P7").Value = 0.2500660443049799
$ws.Range("Q7").Value = 1479.089009612115
$ws.Range("R7").Value = 13311.80108650903
$ws.Range("S7").Value = 0.2283958288653169
$ws.Range("T7").Value = 0.2283958288653169

$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 5.241530333333333
$ws.Range("H8").Value = 15.724591
$ws.Range("I8").Value = 0.04506095322284893
$ws.Range("J8").Value = 0.04506095322284893
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 29.57110033333333
$ws.Range("N8").Value = 88.713301
$ws.Range("O8").Value = 0.5311523066901919
$ws.Range("P8").Value = 0.5311523066901919
$ws.Range("Q8").Value = 154.9978193872101
$ws.Range("R8").Value = 1394.980374484891
$ws.Range("S8").Value = 0.02393422924597505
$ws.Range("T8").Value = 0.02393422924597505

$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 5.241530333333333
$ws.Range("H9").Value = 15.724591
$ws.Range("I9").Value = 0.04506095322284893
$ws.Range("J9").Value = 0.04506095322284893
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 12.180337
$ws.Range("N9").Value = 36.541011
$ws.Range("O9").Value = 0.2187816490048282
$ws.Range("P9").Value = 0.2187816490048282
$ws.Range("Q9").Value = 63.84360585572234
$ws.Range("R9").Value = 574.5924527015011
$ws.Range("S9").Value = 0.009858509651824318
$ws.Range("T9").Value = 0.009858509651824318

$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 5.241530333333333
$ws.Range("H10").Value = 15.724591
$ws.Range("I10").Value = 0.04506095322284893
$ws.Range("J10").Value = 0.04506095322284893
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 13.92204833333333
$ws.Range("N10").Value = 41.76614499999999
$ws.Range("O10").Value = 0.2500660443049799
$ws.Range("P10").Value = 0.2500660443049799
$ws.Range("Q10").Value = 72.97283864129943
$ws.Range("R10").Value = 656.755547771695
$ws.Range("S10").Value = 0.01126821432504957
$ws.Range("T10").Value = 0.01126821432504957
